$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38, shifting existing rows (and table data) down.
$ws.Range("A38:E38").Insert("Down") | Out-Null

# Re-expand the table (ListObject) to cover the newly added row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E133"))

# Populate the new row with the derived-variable data.
$ws.Range("A38").Value = "Ca10j"
$ws.Range("B38").Value = "any_local"
$ws.Range("C38").Value = "Cancer"
$ws.Range("D38").Value = "Any local therapy (surgery or RT) within 3 months"

# Update the view to reflect where the edit happened.
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("D38").Select() | Out-Null
